# New weekly record: insert a row above the current row 141 (Fecha 2021-07-09 /
# serial 44386, "Feria Lagunitas de Puerto Montt" - Cilantro) so every existing
# row from 141 downward shifts down by one (141->142, ..., 171->172), and
# populate the freshly inserted row 141 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 141:171 down to 142:172, leaving row 141 empty for the new record.
$ws.Rows.Item(141).Insert()

# Fill the new row 141 with the new weekly data point.
$ws.Cells.Item(141, 1).Value = 4
$ws.Cells.Item(141, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(141, 3).Value = "Los Lagos"
$ws.Cells.Item(141, 4).Value = 44476
$ws.Cells.Item(141, 5).Value = 10
$ws.Cells.Item(141, 6).Value = 100112040
$ws.Cells.Item(141, 7).Value = "Cilantro"
$ws.Cells.Item(141, 8).Value = "Sin especificar"
$ws.Cells.Item(141, 9).Value = "Primera"
$ws.Cells.Item(141, 10).Value = 200
$ws.Cells.Item(141, 11).Value = 10000
$ws.Cells.Item(141, 12).Value = 10000
$ws.Cells.Item(141, 13).Value = 10000
$ws.Cells.Item(141, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(141, 15).Value = "Región Metropolitana"
$ws.Cells.Item(141, 16).Value = 278
$ws.Cells.Item(141, 17).Value = 36
$ws.Cells.Item(141, 18).Value = "Hortaliza"
